$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.041.04"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = "'1.747.76"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('D4').Value = "'0.9992"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = "'234.38"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.70%  '
$ws.Range('D6').Value = "'0.9992"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = "'0.5295"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.39%  '
$ws.Range('D8').Value = "'0.2788"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.26%  '
$ws.Range('D9').Value = "'0.06188"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.69%  '
$ws.Range('D10').Value = "'1.745.98"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('D11').Value = "'0.07234"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.42%  '
$ws.Range('D12').Value = "'15.37"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').Value = "'0.6446"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('D14').Value = "'4.610"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.72%  '
$ws.Range('D15').Value = "'78.38"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('D16').Value = "'0.9995"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = "'25.960.60"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('E19').Value = '  +1.62%  '
$ws.Range('D20').Value = "'0.000006743"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.97%  '
$ws.Range('D21').Value = "'1.967.95"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'4.317"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.14%  '
$ws.Range('D23').Value = "'8.808"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.74%  '
$ws.Range('D24').Value = "'5.208"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.62%  '
$ws.Range('D25').Value = "'139.11"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('D26').Value = "'1.512"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.69%  '
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('D28').Value = "'1.810"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('D29').Value = "'104.31"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.68%  '
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('D31').Value = "'3.809"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.34%  '
$ws.Range('D32').Value = "'3.684"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.02%  '
$ws.Range('D33').Value = "'0.04541"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.19%  '
$ws.Range('D34').Value = "'2.641"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('D35').Value = "'1.002"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.57%  '
$ws.Range('D36').Value = "'0.6345"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.55%  '
$ws.Range('D37').Value = "'2.712"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.24%  '
$ws.Range('D38').Value = "'0.01597"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.69%  '
$ws.Range('D39').Value = "'1.942"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('D40').Value = "'0.9988"
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Value = "'98.72"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('D42').Value = "'0.3918"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('D43').Value = "'0.7405"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.18%  '
$ws.Range('D44').Value = "'5.038"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.44%  '
$ws.Range('D45').Value = "'0.1146"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.27%  '
$ws.Range('D46').Value = "'6.339"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.62%  '
$ws.Range('D47').Value = "'0.05345"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.50%  '
$ws.Range('D48').Value = "'54.09"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.10%  '
$ws.Range('D49').Value = "'30.71"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.39%  '
$ws.Range('D50').Value = "'7.708"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.50%  '
$ws.Range('D51').Value = "'0.3470"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.50%  '
